# Updated symbol list on Sun Dec 25 16:01:04 UTC 2022 with GitHub Actions
#
# Applies the per-row Price (D) / Volume(1h) label (E) / Hora (G) refresh
# captured in the commit diff. Column G moves from "15" to "16" for every
# data row (2-51); column D gets refreshed price quotes for the rows whose
# price actually moved; two E cells swap their "Bestin24h" / "Worstin24h"
# suffix between KickToken (row 41) and CoinbaseStockToken (row 48).
#
# Cells D/G are plain text (t="inlineStr") in the source workbook, and some
# of the new price strings carry significant trailing zeros (e.g. "5.410"),
# so NumberFormat is forced to "@" (Text) before the assignment -- otherwise
# Excel's COM layer would auto-coerce the numeric-looking string to a
# Number and silently drop the trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G ("Hora"): every data row 2-51 goes from 15 -> 16 ---------
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "16"

# --- Column D ("Price"): refreshed quotes (rows without a price change ---
# --- in the diff are left untouched) --------------------------------------
$priceUpdates = @{
    "D2"  = "243.49"
    "D3"  = "22.89"
    "D4"  = "5.410"
    "D5"  = "0.05929"
    "D7"  = "0.8093"
    "D8"  = "0.9117"
    "D9"  = "0.1418"
    "D10" = "0.07436"
    "D11" = "0.03325"
    "D12" = "0.03086"
    "D13" = "0.09324"
    "D14" = "3.948"
    "D15" = "0.001579"
    "D16" = "0.04802"
    "D18" = "0.005547"
    "D19" = "0.004434"
    "D20" = "0.0009825"
    "D21" = "0.00007810"
    "D22" = "3.653"
    "D23" = "6.437"
    "D40" = "0.03892"
    "D41" = "0.006222"
    "D42" = "0.1067"
    "D43" = "0.003004"
    "D44" = "0.006555"
    "D45" = "0.00005181"
    "D49" = "0.002266"
    "D50" = "0.00002103"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $priceUpdates[$addr]
}

# --- Column E ("Volume(1h)"): Bestin24h badge moves from row 48 to row 41 -
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"
